# This script applies a series of proofing-mark / spell-check-run-split
# edits plus new "results" content to the dissertation draft, matching
# the target OOXML diff. Because w:proofErr marks are not something the
# Word OM exposes a dedicated "add squiggly underline" API for, each
# touched paragraph's Range is rewritten in place via Range.InsertXML
# with hand-built WordprocessingML that reproduces the exact run/proofErr
# structure from the diff while preserving the run text and formatting.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $xml = "<w:p $wNs>$innerXml</w:p>"
    $p.Range.InsertXML($xml)
}

# --- 1. "Hangi data'yi kullanmislar?" (paragraph 3) -------------------
$inner = @"
<w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Hangi</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>data&#8217;yi</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>kullanmislar</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>?</w:t></w:r>
"@
Set-ParaXml 3 $inner

# --- 2. "Bantal et al. ... Messidor databse ..." (paragraph 4) --------
$inner = @"
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Bantal</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> et al. used publicly available MESSIDOR database which includes 1200 images. In this database, there are 3 size of images which are 440 x 960, 2240 x 1488 and 2304 x 1536 pixels. In </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Messidor</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>databse</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> there also information of grading score provided. Every image graded from R0 to R3. With grade R0 means that the patient does not have DR. R1 and R2 are mild and severe cases and if a patient has R3 it means that this is a serious condition of DR. </w:t></w:r>
"@
Set-ParaXml 4 $inner

# --- 3. "Hangi yontemleri uygulamislar?" (paragraph 11) ----------------
$inner = @"
<w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Hangi</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>yontemleri</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>uygulamislar</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>?</w:t></w:r>
"@
Set-ParaXml 11 $inner

# --- 4. "For selecting ensembles several wel-known classifiers ..." (paragraph 25)
$inner = @"
<w:r><w:t xml:space="preserve">For selecting ensembles several </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>wel</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>-known classifiers are trained. These classifiers are</w:t></w:r>
"@
Set-ParaXml 25 $inner

# --- 5. "kNN" (paragraph 27) -------------------------------------------
$inner = @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>kNN</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
"@
Set-ParaXml 27 $inner

# --- 6. "AdaBoost" (paragraph 28) --------------------------------------
$inner = @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>AdaBoost</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
"@
Set-ParaXml 28 $inner

# --- 7. "X2-X7 represent the result of MA detection. Xi, i represents..." (paragraph 40)
$inner = @"
<w:r><w:t xml:space="preserve">X2-X7 represent the result of MA detection. Xi, </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>i</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> represents number of MAs at the confidence levels&#8230;&#8230;.</w:t></w:r>
"@
Set-ParaXml 40 $inner

# --- 8. "Neyle test etmisler? Sonuclari neler?" (paragraph 46) --------
$inner = @"
<w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Neyle</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> test </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>etmisler</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>?</w:t></w:r>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>S</w:t></w:r>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>onuclari</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>neler</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>?</w:t></w:r>
"@
Set-ParaXml 46 $inner

# --- 9. "To compare their results with others ... Reciever ..." (paragraph 52)
$inner = @"
<w:r><w:t xml:space="preserve">To compare their results with others they have fitted </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Reciever</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> Operating Characteristic curves to the results and calculated AUC using JROCFIR</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
"@
Set-ParaXml 52 $inner

# --- 10-13. Tail: "R0 vs {R1, R2, R3}" + bookmark paragraph (57-58) ---
# replaced & extended with the new results narrative.
$p57 = $d.Paragraphs(57)
$p58 = $d.Paragraphs(58)
$tailRange = $d.Range($p57.Range.Start, $p58.Range.End)
$tailXml = @"
<w:p $wNs>
  <w:r><w:t>R</w:t></w:r>
  <w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>0</w:t></w:r>
  <w:r><w:t xml:space="preserve"> vs {R</w:t></w:r>
  <w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>1</w:t></w:r>
  <w:r><w:t>, R</w:t></w:r>
  <w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>2</w:t></w:r>
  <w:r><w:t>, R</w:t></w:r>
  <w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>3</w:t></w:r>
  <w:r><w:t>} = No DR vs DR</w:t></w:r>
</w:p>
<w:p $wNs/>
<w:p $wNs>
  <w:pPr><w:rPr><w:i/></w:rPr></w:pPr>
  <w:r><w:t xml:space="preserve">For R0 vs R1 best performing ensemble achieved 94% Sensitivity, 90% Accuracy, 90% Specificity using backward search, output fusion strategy </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>D</w:t></w:r>
  <w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>avg</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> and energy function </w:t></w:r>
  <w:r><w:rPr><w:i/></w:rPr><w:t>Accuracy</w:t></w:r>
  <w:r><w:rPr><w:i/></w:rPr><w:t>.</w:t></w:r>
</w:p>
<w:p $wNs/>
<w:p $wNs/>
<w:p $wNs>
  <w:r><w:t xml:space="preserve">No DR vs DR; 90% Sensitivity, 91% Specificity and 90% Accuracy with same search method and fusion strategy but different energy </w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t xml:space="preserve">function which is Sensitivity. </w:t></w:r>
</w:p>
"@
$tailRange.InsertXML($tailXml)
